$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Chlorophyll A (92nd Percentile)" record for the 2017-2021 period
# (row 7) was removed; the rows below it (MCI, QMCI for 2017-2021) shift
# up to take its place.
$ws.Rows.Item(7).Delete()
